# Update the "dSF" column (F) values for the specified rows, per the
# repull/recalculation of data described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = -8
$ws.Range("F3").Value  = -3
$ws.Range("F5").Value  = -3
$ws.Range("F6").Value  = -1
$ws.Range("F10").Value = -4
$ws.Range("F11").Value = -17
$ws.Range("F14").Value = -4
$ws.Range("F16").Value = -3
$ws.Range("F21").Value = -4
